# Commit: "Fruta / hortaliza, semanal"
# Insert 3 new weekly price rows for "Hass" (1a/2a/3a nueva(o)) dated 44489
# right after the current row 174, pushing the existing rows 175-253 down to
# 178-256 (the sheet's dimension grows from A1:T253 to A1:T256).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 175; this shifts old rows 175:253 -> 178:256
$ws.Rows("175:177").Insert()

# Constant column values shared by every data row in this block
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$unidad17  = "`$/kilo (en caja de 17 kilos)"
$origen    = "Provincia de Limarí"
$kgUnidad  = 1

$fecha = 44489
$variedad = "Hass"

# New row 175: Hass / 1a nueva(o)
$r = 175
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "1a nueva(o)"
$ws.Cells.Item($r, 13).Value = 400
$ws.Cells.Item($r, 14).Value = 2400
$ws.Cells.Item($r, 15).Value = 2500
$ws.Cells.Item($r, 16).Value = 2450
$ws.Cells.Item($r, 17).Value = $unidad17
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2450
$ws.Cells.Item($r, 20).Value = $kgUnidad

# New row 176: Hass / 2a nueva(o)
$r = 176
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "2a nueva(o)"
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 2100
$ws.Cells.Item($r, 15).Value = 2200
$ws.Cells.Item($r, 16).Value = 2150
$ws.Cells.Item($r, 17).Value = $unidad17
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2150
$ws.Cells.Item($r, 20).Value = $kgUnidad

# New row 177: Hass / 3a nueva (o)
$r = 177
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "3a nueva (o)"
$ws.Cells.Item($r, 13).Value = 240
$ws.Cells.Item($r, 14).Value = 1800
$ws.Cells.Item($r, 15).Value = 1900
$ws.Cells.Item($r, 16).Value = 1850
$ws.Cells.Item($r, 17).Value = $unidad17
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1850
$ws.Cells.Item($r, 20).Value = $kgUnidad

Write-Output "Inserted rows 175-177; dimension now $($ws.UsedRange.Address())"
